# Book excel half update
# Converts the small NUM1/NUM2 demo sheet into a Library Management
# "Book" entry sheet: new headers, one data row, header styling
# (bold Arial font, a filled/centred pair of header cells), and
# print/page setup + selection tweaks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Get rid of the old demo rows/values we no longer need.
# ------------------------------------------------------------------
# Row 3 (A3=4, B3=5) is not part of the new layout at all.
$ws.Rows.Item(3).Delete() | Out-Null
# A2 stays empty in the new layout (row 2 data starts at column B).
$ws.Range("A2").ClearContents() | Out-Null

# ------------------------------------------------------------------
# 2. New header row (row 1) + data row (row 2).
#    Values are entered in the same order the original author typed
#    them (plain headers/data first, the "Title 3" header in E1
#    added last) so the shared-string table comes out in the same
#    order as the authored workbook.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Book No"
$ws.Range("C1").Value = "Author 1"
$ws.Range("D1").Value = "Author 2"
$ws.Range("B1").Value = "Title"
$ws.Range("F1").Value = "Edition"
$ws.Range("G1").Value = "Publisher"
$ws.Range("H1").Value = "CL Number"
$ws.Range("I1").Value = "Total Pages"
$ws.Range("J1").Value = "Cost"
$ws.Range("K1").Value = "Supplier"
$ws.Range("L1").Value = "Remark"
$ws.Range("M1").Value = "Bill_No"
# N1 carries the header styling too, but no text.

# Row 2 - a single sample book record.
$ws.Range("B2").Value = "holly cow"
$ws.Range("C2").Value = "hijfds"
$ws.Range("D2").Value = "hghvfhjg"
$ws.Range("E2").Value = "hjkgkgvhjgk"
$ws.Range("F2").Value = "hj,vbhjkg"
$ws.Range("G2").Value = "hv,bjhjk"
$ws.Range("H2").Value = 6.4
$ws.Range("I2").Value = 97
$ws.Range("J2").Value = 351
$ws.Range("K2").Value = "jkkhg"
$ws.Range("L2").Value = "hjkgkhjgvgj"
$ws.Range("M2").Value = "jkhv"

# "Title 3" header, added last.
$ws.Range("E1").Value = "Title 3"

# ------------------------------------------------------------------
# 4. Header formatting.
#    Most header cells: bold 8pt Arial, black text, no fill.
#    Publisher/CL Number (G1:H1): same bold font plus a white fill
#    and left/centre aligned text.
# ------------------------------------------------------------------
# (Comma-joined / Union() ranges aren't reliable in this host, so the
# two contiguous blocks that share formatting are styled separately.)
$plainHeaders1 = $ws.Range("A1:F1")
$plainHeaders1.Font.Name = "Arial"
$plainHeaders1.Font.Size = 8
$plainHeaders1.Font.Bold = $true
$plainHeaders1.Font.Color = 0

$plainHeaders2 = $ws.Range("I1:N1")
$plainHeaders2.Font.Name = "Arial"
$plainHeaders2.Font.Size = 8
$plainHeaders2.Font.Bold = $true
$plainHeaders2.Font.Color = 0

$fancyHeaders = $ws.Range("G1:H1")
$fancyHeaders.Font.Name = "Arial"
$fancyHeaders.Font.Size = 8
$fancyHeaders.Font.Bold = $true
$fancyHeaders.Font.Color = 0
$fancyHeaders.Interior.Color = 16777215
$fancyHeaders.HorizontalAlignment = -4131
$fancyHeaders.VerticalAlignment = -4108

# ------------------------------------------------------------------
# 5. Selection + page setup.
# ------------------------------------------------------------------
$ws.Range("G9").Select() | Out-Null

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
